$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table on the "CLX + Morpheus" slide: switch built-in table style
#    from "No Style, Table Grid" to "No Style, No Grid".
# ---------------------------------------------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle("{22E53491-32EE-45D4-8380-BA633CB99490}", $true)
        }
    }
}

# ---------------------------------------------------------------------
# 2) Last slide: the "@RAPIDSai" twitter-handle run should no longer use
#    the theme hyperlink color - it should match the surrounding white
#    (lt1) text color.
# ---------------------------------------------------------------------
$lastSlide = $p.Slides.Item($p.Slides.Count)
for ($i = 1; $i -le $lastSlide.Shapes.Count; $i++) {
    $shp = $lastSlide.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        $runCount = $tr.Runs().Count
        for ($j = 1; $j -le $runCount; $j++) {
            $run = $tr.Runs($j)
            if ($run.Text.Trim() -eq "@RAPIDSai") {
                $run.Font.Color.ObjectThemeColor = 2   # msoThemeColorLight1 -> schemeClr val="lt1"
            }
        }
    }
}
